# Witval stage 3 evidence
# Fill in the newly-collected TxHash evidence for sheets B1, B2, B5, B6
# and add a new sheet B7 with its own evidence, mirroring the upstream commit.

$wb = $excel.ActiveWorkbook

# --- B1 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("B1")
$ws.Range("A2").Value = "519FF890B2C8A68C347B9DA7D6DF2E59D9248A04852326EED293408CCDD59547"
$ws.Range("A3").Value = "A481851C99D767E88C30F96C61C48E0FB1FEB80A84AA029F55985CD894C3FF62"
$ws.Range("A3").Select() | Out-Null

# --- B2 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("B2")
$ws.Range("A2").Value = "C808D28FDE0BEE6864ADF8142B75B4B3BFABCDE5A00F1F5A249E46242BFFCC8A"
$ws.Range("A3").Value = "5992C62F45ADECA548909285CF400FDDB82ACA3F75C0D05264CFE26FF279FBE3"
$ws.Range("A3").Select() | Out-Null

# --- B5 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("B5")
$ws.Range("A2").Value = "7AC8546B236E4672FABF22AE411489F2A5308E6CF194934E616E13C3BD99B82C"
$ws.Range("A3").Value = "91B12BBCCE8ECD7CC508BA7E76643452018DCAFDC8ABA3E185453B5A19414CC6"
$ws.Range("A6").Select() | Out-Null

# --- B6 -------------------------------------------------------------
$ws = $wb.Worksheets.Item("B6")
$ws.Range("A2").Value = "BA1C1DD06721CC4C843181F1A4745B40599B180B26A856F54FE2976EFD4CC2E0"
$ws.Range("A3").Value = "7987FA76335C203310E2D27C4C0A859B95B52FE75F962E12E4ED3082F785086A"
$ws.Range("A1").Select() | Out-Null

# --- B7 (new sheet, added after B6) ---------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "B7"
$newSheet.Range("A1").Value = "TxHash"
$newSheet.Range("A2").Value = "1E85392C414980C7AE669C7AB358D03D238B47F8480991BCA667BD80D17C5835"
$newSheet.Range("A3").Value = "BC6A74F171548B7CDF8F3B74058BDEBABB07FFDE5A52343A3AB75D16BAEFCF94"

# Match the header styling used by the other evidence sheets (e.g. B6!A1)
$srcSheet = $wb.Worksheets.Item("B6")
$srcSheet.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("A3").Select() | Out-Null
